$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The capacitor rows (25-28) are being re-shuffled:
#   old row25 "C3"    (1.0 uF cap)   -> moves down, becomes new row28 "C9"
#   old row26 "C4-C6" (0.1 uF cap x3)-> moves up,   becomes new row25 "C3-C6"
#   old row27 "C7"    (0.33 uF cap)  -> moves up,   becomes new row26 "C7" (unchanged name)
#   old row28 "C8"    (leave-out note)-> stays,     becomes new row27 "C8" (unchanged name),
#                                        but loses every column value except A/B
# ---------------------------------------------------------------------------

# --- Row 28 (becomes "C9", gets the old row25 "C3" content) -------------------
$ws.Range("A28").Value = "C9"
$ws.Range("B28").Value = "1.0 uF cap"
$ws.Range("C28").Value = "GRJ31MR71C105KE01L"
$ws.Range("D28").Value = "Murata"
$ws.Range("E28").Value = "C1206"
$ws.Range("F28").Value = "Mouser"
$ws.Range("G28").Value = "81-GRJ31MR71C105KE1L"
$ws.Range("H28").Value = 0.10199999999999999
$ws.Range("I28").Value = 1
$ws.Range("J28").Formula = "=I28*H28"

# --- Row 25 (becomes "C3-C6", gets the old row26 "C4-C6" content) -------------
$ws.Range("A25").Value = "C3-C6"
$ws.Range("B25").Value = "0.1 uF cap"
$ws.Range("C25").Value = "GRM21BR71H104KA01L"
$ws.Range("D25").Value = "Murata"
$ws.Range("E25").Value = "C0805"
$ws.Range("F25").Value = "Mouser"
$ws.Range("G25").Value = "81-GRM40X104K50L"
$ws.Range("H25").Value = 0.042000000000000003
$ws.Range("I25").Value = 3

# --- Row 26 (becomes "C7", gets the old row27 "C7" content) -------------------
$ws.Range("A26").Value = "C7"
$ws.Range("B26").Value = "0.33 uF cap"
$ws.Range("C26").Value = "GRM21BR71C334KA01L"
$ws.Range("D26").Value = "Murata"
$ws.Range("E26").Value = "C1206"
$ws.Range("F26").Value = "Mouser"
$ws.Range("G26").Value = "81-GRM40X334K16L"
$ws.Range("H26").Value = 0.069000000000000006
$ws.Range("I26").Value = 1

# C26 picks up the font style that used to live on C27 (and vice-versa below)
$ws.Range("Z1").Value = ""
$ws.Range("C27").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("C26").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# --- Row 27 (becomes "C8", keeps only the note, loses the rest) ---------------
$ws.Range("A27").Value = "C8"
$ws.Range("B27").Value = "leave out, better to put on motor jack"
$ws.Range("C27").ClearContents()
$ws.Range("D27").Clear()
$ws.Range("E27").Clear()
$ws.Range("F27").Clear()
$ws.Range("G27").ClearContents()
$ws.Range("H27").Clear()
$ws.Range("I27").ClearContents()
$ws.Range("J27").Clear()

# --- Totals: the sum now needs to cover the new row28 as well -----------------
$ws.Range("J31").Formula = "=SUM(J2:J28)"

# --- Sheet / page setup tweaks --------------------------------------------
$ws.PageSetup.Orientation = 2
$ws.PageSetup.Zoom = $false
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = $false

$ws.Range("B34").Select()
